$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the bordered header style (currently on A1, style index 1) into a
# holding cell well outside the diagram so it survives the upcoming clear.
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Wipe the old diagram (A1:G11) completely.
$ws.Range("A1:G11").Clear()

# Stamp the bordered style onto the two new header rows from the holding
# cell, so the cells land on the *existing* style slot instead of Excel
# minting a fresh (visually identical) one.
$ws.Range("Z1").Copy()
$ws.Range("B5:H5").PasteSpecial(-4122)
$ws.Range("B10:G10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("Z1").Clear()

# --- First frame diagram (rows 5-7, cols B-H) ---
$ws.Range("E5").Value = "frameNo 1 byte"
$ws.Range("F5").Value = "data 2 bytes"
$ws.Range("G5").Value = "check 1 byte"
$ws.Range("H5").Value = "endFlag 1 byte"
$ws.Range("D5").Value = "src 1 byte"
$ws.Range("C5").Value = "dst 1 byte"
$ws.Range("B5").Value = "openFlag 1 byte"

$ws.Range("B6").Value = "char"
$ws.Range("C6").Value = "char"
$ws.Range("D6").Value = "char"
$ws.Range("E6").Value = "char"
$ws.Range("F6").Value = "char"
$ws.Range("G6").Value = "char"
$ws.Range("H6").Value = "char end 1 more"

$ws.Range("H7").Value = "char end 0 end"

# --- Second frame diagram (rows 10-12, cols B-G) ---
$ws.Range("E10").Value = "frameNo 1 byte"
$ws.Range("F10").Value = "check 1 byte"
$ws.Range("G10").Value = "endFlag 1 byte"
$ws.Range("D10").Value = "src 1 byte"
$ws.Range("C10").Value = "dst 1 byte"
$ws.Range("B10").Value = "openFlag 1 byte"

$ws.Range("B11").Value = "char"
$ws.Range("C11").Value = "char"
$ws.Range("D11").Value = "char"
$ws.Range("E11").Value = "char"
$ws.Range("F11").Value = "char"
$ws.Range("G11").Value = "char end 1 more"

$ws.Range("G12").Value = "char end 0 end"

# --- Label ---
$ws.Range("B15").Value = "stop and wait arq"

# Auto-fit columns to content, matching the author's resize.
$ws.Columns("B:H").AutoFit()

$ws.Range("H21").Select()
